$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.746.70'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '1.856.16'
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.021'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -1.54%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.018'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4377'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3788'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07435'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8836'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.55'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("D12").Value = '1.846.80'
$ws.Range("E12").Value = '  -1.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.792'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.494'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07143'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("E16").Value = '  +5.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.021'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009027'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.018'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("D21").Value = '27.739.09'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.277'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.50%  '
$ws.Range("D24").Value = '2.080.15'
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.033'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.71'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.449'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.991'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.87'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09040'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.230'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7689'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.013'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.565'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("E36").Value = '  -1.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.141'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01981'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05306'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.862'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5188'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.962'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1678'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.711'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '110.05'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.52%  '
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4734'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.020'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06470'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.847'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.91%  '
